$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update individual cell values per the diff
$ws.Range("C3").Value = -12.912
$ws.Range("B9").Value = 6.484999999999999
$ws.Range("B18").Value = 5.972
$ws.Range("B20").Value = 6.37
$ws.Range("D21").Value = -7.805
